# Entered entry with UserId:104
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# UserId in row 2 was corrected from 100 to 103
$ws.Range("A2").Value = 103

# New row 3: UserId 104, Agricultural Irrigation / Article entry
$ws.Cells.Item(3, 1).Value = 104
$ws.Cells.Item(3, 2).Value = "Agricultural Irrigation"
$ws.Cells.Item(3, 3).Value = 9
$ws.Cells.Item(3, 4).Value = "Article"
$ws.Cells.Item(3, 5).Value = "435-235-656-232"
# Publisher entered before Author to match original authoring order
$ws.Cells.Item(3, 7).Value = "Covenant Uni. Research"
$ws.Cells.Item(3, 6).Value = "Covenant University"

$ws.Range("I3").Select() | Out-Null
